# Loan RBI, Variable Instalments
#
# On the "Repayment Schedule" sheet, insert a new blank column before the
# existing "N" column (the old "Late" column and everything to its right
# shifts one column over: N->O, O->P, P->Q). The new column N is left
# blank (header + data) but is given an explicit width of 10, matching
# the width of the "In Advance" column next to it.
#
# The workbook's active sheet also moves from "NewLoanInput" to
# "Repayment Schedule", with the selected cell there becoming T7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment Schedule")

# Insert a new blank column at N (shifts N:P -> O:Q).
$ws.Columns("N:N").Insert() | Out-Null

# Give the newly inserted column an explicit (non-autofit) width of 10.
$ws.Columns("N:N").ColumnWidth = 9.17

# Make "Repayment Schedule" the active sheet and select T7 on it
# (this clears tabSelected on whichever sheet was active before,
# e.g. "NewLoanInput", and sets it here instead).
$ws.Activate() | Out-Null
$ws.Range("T7").Select() | Out-Null
